# Daily auto push: append the 2025/09/29 17:00 row (rank 3) to the tracking
# sheet, extending the used range from A1:D34 to A1:D35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (matching the existing rows, which
# are literal "yyyy/mm/dd" strings rather than real date values). Force the
# cell to Text format before assigning so Excel doesn't auto-convert the
# "2025/09/29" string into a date serial, then clear the format back off so
# the cell keeps the sheet's default (unstyled) look, same as every other
# data row.
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "2025/09/29"
$ws.Range("A35").ClearFormats()

$ws.Range("B35").Value = "月"
$ws.Range("C35").Value = 17
$ws.Range("D35").Value = 3
